$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Lrp6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.07131444737854614
$ws.Range("J2").Value = 0.07131444737854616
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.210619666666666
$ws.Range("N2").Value = 27.631859
$ws.Range("O2").Value = 0.133636377806767
$ws.Range("P2").Value = 0.133636377806767
$ws.Range("Q2").Value = 0.249819637219
$ws.Range("R2").Value = 2.248376734971
$ws.Range("S2").Value = 0.009530204432960196
$ws.Range("T2").Value = 0.009530204432960198

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Lrp6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.07131444737854614
$ws.Range("J3").Value = 0.07131444737854616
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.97212233333332
$ws.Range("N3").Value = 131.916367
$ws.Range("O3").Value = 0.6379891218794987
$ws.Range("P3").Value = 0.6379891218794989
$ws.Range("Q3").Value = 1.192655874047
$ws.Range("R3").Value = 10.733902866423
$ws.Range("S3").Value = 0.04549784166036037
$ws.Range("T3").Value = 0.04549784166036039

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Lrp6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.07131444737854614
$ws.Range("J4").Value = 0.07131444737854616
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.740255
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2283745003137342
$ws.Range("P4").Value = 0.2283745003137342
$ws.Range("Q4").Value = 0.4269229363649999
$ws.Range("R4").Value = 3.842306427285
$ws.Range("S4").Value = 0.01628640128522556
$ws.Range("T4").Value = 0.01628640128522557

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Lrp6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3532066666666667
$ws.Range("H5").Value = 1.05962
$ws.Range("I5").Value = 0.9286855526214538
$ws.Range("J5").Value = 0.9286855526214538
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.210619666666666
$ws.Range("N5").Value = 27.631859
$ws.Range("O5").Value = 0.133636377806767
$ws.Range("P5").Value = 0.133636377806767
$ws.Range("Q5").Value = 3.253252270397778
$ws.Range("R5").Value = 29.27927043358
$ws.Range("S5").Value = 0.1241061733738068
$ws.Range("T5").Value = 0.1241061733738068

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Lrp6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3532066666666667
$ws.Range("H6").Value = 1.05962
$ws.Range("I6").Value = 0.9286855526214538
$ws.Range("J6").Value = 0.9286855526214538
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 43.97212233333332
$ws.Range("N6").Value = 131.916367
$ws.Range("O6").Value = 0.6379891218794987
$ws.Range("P6").Value = 0.6379891218794989
$ws.Range("Q6").Value = 15.53124675561555
$ws.Range("R6").Value = 139.78122080054
$ws.Range("S6").Value = 0.5924912802191383
$ws.Range("T6").Value = 0.5924912802191384

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Lrp6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3532066666666667
$ws.Range("H7").Value = 1.05962
$ws.Range("I7").Value = 0.9286855526214538
$ws.Range("J7").Value = 0.9286855526214538
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.740255
$ws.Range("N7").Value = 47.220765
$ws.Range("O7").Value = 0.2283745003137342
$ws.Range("P7").Value = 0.2283745003137342
$ws.Range("Q7").Value = 5.559563001033333
$ws.Range("R7").Value = 50.0360670093
$ws.Range("S7").Value = 0.2120880990285086
$ws.Range("T7").Value = 0.2120880990285086
